$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Collin Sexton -> Paul George
$ws.Range("A3").Value = "Paul George"
$ws.Range("B3").Value = "SG,SF,PF"
$ws.Range("C3").Value = "Philadelphia 76ers"

# Row 6: Caris LeVert -> Alex Caruso (position unchanged)
$ws.Range("A6").Value = "Alex Caruso"
$ws.Range("C6").Value = "Oklahoma City Thunder"

# Row 7: Jayson Tatum -> Draymond Green
$ws.Range("A7").Value = "Draymond Green"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Golden State Warriors"

# Row 14: Paul George -> Collin Sexton
$ws.Range("A14").Value = "Collin Sexton"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Utah Jazz"

# Row 15: Draymond Green -> Jayson Tatum
$ws.Range("A15").Value = "Jayson Tatum"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Boston Celtics"
